# Odin_Profits workbook refresh (scheduled runner).
# Leve crafting-profit figures (currentAveragePrice / NQ / HQ, LevePrice NQ/HQ,
# LeveProfit NQ/HQ) were recomputed from a new market-board price pull.
# All edits below are plain value writes -- the source sheets contain no formulas.

$wb = $excel.ActiveWorkbook

# ALC!row 86
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 3166.6667
$ws.Range("I86").Value = 3500
$ws.Range("J86").Value = 2500
$ws.Range("K86").Value = 3500
$ws.Range("L86").Value = 2500
$ws.Range("M86").Value = -2377
$ws.Range("N86").Value = -4746

# ALC!row 89
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 3166.6667
$ws.Range("I89").Value = 3500
$ws.Range("J89").Value = 2500
$ws.Range("K89").Value = 17500
$ws.Range("L89").Value = 12500
$ws.Range("M89").Value = -11884
$ws.Range("N89").Value = -23732

# ALC!row 113
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 52002.5
$ws.Range("I113").Value = 4005
$ws.Range("J113").Value = 100000
$ws.Range("K113").Value = 4005
$ws.Range("L113").Value = 100000
$ws.Range("M113").Value = -751
$ws.Range("N113").Value = -106508

# ALC!row 116
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 29999
$ws.Range("J116").Value = 19998
$ws.Range("L116").Value = 19998
$ws.Range("N116").Value = -26882

# ARM!row 2
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 13943.556
$ws.Range("I2").Value = 7166.6665
$ws.Range("K2").Value = 7166.6665
$ws.Range("M2").Value = -7053.6665

# ARM!row 45
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 5240
$ws.Range("I45").Value = 2800
$ws.Range("K45").Value = 2800
$ws.Range("M45").Value = -2423

# ARM!row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 6393.8276
$ws.Range("I61").Value = 7225.615
$ws.Range("K61").Value = 7225.615
$ws.Range("M61").Value = -7013.615

# ARM!row 116
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 13943.556
$ws.Range("I116").Value = 7166.6665
$ws.Range("K116").Value = 7166.6665
$ws.Range("M116").Value = -4872.6665

# ARM!row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 6393.8276
$ws.Range("I136").Value = 7225.615
$ws.Range("K136").Value = 21676.845
$ws.Range("M136").Value = -19126.845

# BSM!row 3
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 13943.556
$ws.Range("I3").Value = 7166.6665
$ws.Range("K3").Value = 7166.6665
$ws.Range("M3").Value = -7052.6665

# CRP!row 58
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 13338.571
$ws.Range("I58").Value = 7973.6665
$ws.Range("J58").Value = 22995.4
$ws.Range("K58").Value = 7973.6665
$ws.Range("L58").Value = 22995.4
$ws.Range("M58").Value = -7770.6665
$ws.Range("N58").Value = -23401.4

# CRP!row 86
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 21092.223
$ws.Range("J86").Value = 18899.5
$ws.Range("L86").Value = 18899.5
$ws.Range("N86").Value = -21145.5

# CRP!row 89
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H89").Value = 21092.223
$ws.Range("J89").Value = 18899.5
$ws.Range("L89").Value = 94497.5
$ws.Range("N89").Value = -105729.5

# CRP!row 99
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 7888
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 7888
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 7888
$ws.Range("M99").ClearContents()  # was -1624
$ws.Range("N99").Value = -10884

# CRP!row 107
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 2999
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 2999
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 2999
$ws.Range("M107").ClearContents()  # was 1621
$ws.Range("N107").Value = -6839

# CRP!row 126
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 7888
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 7888
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 23664
$ws.Range("M126").ClearContents()  # was -6896
$ws.Range("N126").Value = -28604

# CRP!row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 15783.875
$ws.Range("I132").Value = 28954
$ws.Range("J132").Value = 11393.833
$ws.Range("K132").Value = 86862
$ws.Range("L132").Value = 34181.499
$ws.Range("M132").Value = -84332
$ws.Range("N132").Value = -39241.499

# CRP!row 136
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 13338.571
$ws.Range("I136").Value = 7973.6665
$ws.Range("J136").Value = 22995.4
$ws.Range("K136").Value = 23920.9995
$ws.Range("L136").Value = 68986.20000000001
$ws.Range("M136").Value = -21370.9995
$ws.Range("N136").Value = -74086.20000000001

# CUL!row 14
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 236.88889
$ws.Range("I14").Value = 236.88889
$ws.Range("K14").Value = 710.6666700000001
$ws.Range("M14").Value = -537.6666700000001

# CUL!row 39
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 5160
$ws.Range("J39").Value = 5160
$ws.Range("L39").Value = 15480
$ws.Range("N39").Value = -16068

# CUL!row 62
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H62").Value = 13714
$ws.Range("J62").Value = 16799.8
$ws.Range("L62").Value = 50399.39999999999
$ws.Range("N62").Value = -51771.39999999999

# CUL!row 65
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H65").Value = 13714
$ws.Range("J65").Value = 16799.8
$ws.Range("L65").Value = 151198.2
$ws.Range("N65").Value = -158062.2

# CUL!row 134
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H134").Value = 6692.231
$ws.Range("I134").Value = 5697.875
$ws.Range("K134").Value = 17093.625
$ws.Range("M134").Value = -12023.625

# GSM!row 122
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 4859.4736
$ws.Range("I122").Value = 2258.3
$ws.Range("J122").Value = 7749.6665
$ws.Range("K122").Value = 6774.900000000001
$ws.Range("L122").Value = 23248.9995
$ws.Range("M122").Value = -4324.900000000001
$ws.Range("N122").Value = -28148.9995

# GSM!row 126
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 13113.777
$ws.Range("I126").Value = 2012
$ws.Range("J126").Value = 14501.5
$ws.Range("K126").Value = 6036
$ws.Range("L126").Value = 43504.5
$ws.Range("M126").Value = -3566
$ws.Range("N126").Value = -48444.5

# GSM!row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 5597.3
$ws.Range("I132").Value = 7233.3335
$ws.Range("J132").Value = 4896.143
$ws.Range("K132").Value = 21700.0005
$ws.Range("L132").Value = 14688.429
$ws.Range("M132").Value = -19170.0005
$ws.Range("N132").Value = -19748.429

# LTW!row 22
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1499.6666
$ws.Range("J22").Value = 1999.5
$ws.Range("L22").Value = 1999.5
$ws.Range("N22").Value = -2589.5

# LTW!row 27
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 1499.6666
$ws.Range("J27").Value = 1999.5
$ws.Range("L27").Value = 1999.5
$ws.Range("N27").Value = -2213.5

# LTW!row 46
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1683.25
$ws.Range("I46").Value = 850
$ws.Range("J46").Value = 1802.2858
$ws.Range("K46").Value = 850
$ws.Range("L46").Value = 1802.2858
$ws.Range("M46").Value = -662
$ws.Range("N46").Value = -2178.2858

# LTW!row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2667.861
$ws.Range("I132").Value = 2400.2222
$ws.Range("J132").Value = 3470.7778
$ws.Range("K132").Value = 7200.6666
$ws.Range("L132").Value = 10412.3334
$ws.Range("M132").Value = -4670.6666
$ws.Range("N132").Value = -15472.3334

# LTW!row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 5633.857
$ws.Range("I136").Value = 5463.091
$ws.Range("K136").Value = 16389.273
$ws.Range("M136").Value = -13839.273

# WVR!row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 57197.85
$ws.Range("I136").Value = 64200.5
$ws.Range("K136").Value = 192601.5
$ws.Range("M136").Value = -190051.5
